# Correct thrown exception name
# Slide 1, shape 2 (body text), paragraph that explains the predefined-threshold behaviour.
#
# Before:
#   "If the predefined threshold is exceeded, then the strategy will throw "
#   "RateLimitRejectedException"
#   "."
#
# After:
#   "If the predefined threshold is exceeded, "
#   "then "
#   "it"
#   " will throw an "
#   "RateLimiterRejectedException"
#   "."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$oldRun1 = "If the predefined threshold is exceeded, then the strategy will throw "
$oldRun2 = "RateLimitRejectedException"

$newRun1 = "If the predefined threshold is exceeded, "
$thenText = "then "
$itText = "it"
$throwText = " will throw an "
$newExcText = "RateLimiterRejectedException"

# 1) Shrink the first run down to its new (shorter) text - this keeps the
#    run's original rPr (formatting) completely untouched.
$fullText = $tr.Text
$run1Start = $fullText.IndexOf($oldRun1) + 1
$run1Range = $tr.Characters($run1Start, $oldRun1.Length)
$run1Range.Text = $newRun1

# 2) Insert the replacement wording right after it (still same run/
#    formatting as run 1 for now - it gets split into separate runs with
#    the correct formatting below).
$run1RangeShrunk = $tr.Characters($run1Start, $newRun1.Length)
$run1RangeShrunk.InsertAfter($thenText + $itText + $throwText) | Out-Null

# 3) "then " -> its own (non-bold) run.
$thenStart = $run1Start + $newRun1.Length
$thenRange = $tr.Characters($thenStart, $thenText.Length)
$thenRange.Font.Bold = $false

# 4) "it" -> its own (non-bold) run.
$itStart = $thenStart + $thenText.Length
$itRange = $tr.Characters($itStart, $itText.Length)
$itRange.Font.Bold = $false

# 5) " will throw an " -> its own (non-bold) run.
$throwStart = $itStart + $itText.Length
$throwRange = $tr.Characters($throwStart, $throwText.Length)
$throwRange.Font.Bold = $false

# 6) Rename the exception type in place - keeps the existing run's rPr
#    (dirty/err spell-check flags) untouched, only the text content changes.
$excStart = $throwStart + $throwText.Length
$excRange = $tr.Characters($excStart, $oldRun2.Length)
$excRange.Text = $newExcText
